$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsOverview.Range("G2").Value = "2016-09-06 19:28:47"

$wsZhCn.Range("H2").Value = "2016-09-06 19:28:42"
$wsZhCn.Range("K2").Value = "2016-09-06 19:29:00"

$wsDeDe.Range("K2").Value = "2016-09-06 19:29:18"
